# Rename the spectral-value channel headers from the raw ShakeMap codes
# (pga, pgv, psa03, psa10, psa30) to their display labels
# (PGA, PGV, SA(0.3), SA(1.0), SA(3.0)).  These labels appear three times
# on row 3 of the sheet - once for each of the three channel blocks
# (I:M, N:R, S:W).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @("PGA", "PGV", "SA(0.3)", "SA(1.0)", "SA(3.0)")

$block1 = @("I3", "J3", "K3", "L3", "M3")
$block2 = @("N3", "O3", "P3", "Q3", "R3")
$block3 = @("S3", "T3", "U3", "V3", "W3")

foreach ($i in 0..4) {
    $ws.Range($block1[$i]).Value = $labels[$i]
    $ws.Range($block2[$i]).Value = $labels[$i]
    $ws.Range($block3[$i]).Value = $labels[$i]
}

# Update the current selection to match the author's last-saved cursor
# position (cell L4) instead of the old D2.
$ws.Range("L4").Select()
